$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
try {
  $s.Shapes.Item(1).ThisMethodDoesNotExistAtAll(1,2,3)
  Write-Host "no error"
} catch {
  Write-Host "ERR: $_"
}
